# Fruta / hortaliza, semanal
#
# Inserts one new daily price record for "Nectarín" (Feria Lagunitas de
# Puerto Montt) into the data table. The new record belongs right after
# the current header/first block, at row 397, which pushes every
# existing record from row 397 down by one row (397->398, ..., 511->512)
# and grows the used range from A1:T511 to A1:T512.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 397 (and everything below it) down by one row.
$ws.Rows.Item(397).Insert()

# Populate the newly-opened row 397 with the new record. Columns A, B,
# C, E, F, G, H, I, J are constant for every row in this sheet.
$ws.Range("A397").Value2 = 4
$ws.Range("B397").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C397").Value2 = "Los Lagos"
$ws.Range("D397").Value2 = 44900
$ws.Range("E397").Value2 = 10
$ws.Range("F397").Value2 = "Fruta"
$ws.Range("G397").Value2 = 100103
$ws.Range("H397").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I397").Value2 = 100103006
$ws.Range("J397").Value2 = "Nectarín"
$ws.Range("K397").Value2 = "Early Glo"
$ws.Range("L397").Value2 = "Primera"
$ws.Range("M397").Value2 = 600
$ws.Range("N397").Value2 = 23000
$ws.Range("O397").Value2 = 24000
$ws.Range("P397").Value2 = 23500
$ws.Range("Q397").Value2 = "$/caja 14 kilos empedrada"
$ws.Range("R397").Value2 = "Provincia de San Felipe de Aconcagua"
$ws.Range("S397").Value2 = 1679
$ws.Range("T397").Value2 = 14

# Keep the date column's display format consistent with the rest of the
# column (YYYY-MM-DD HH:MM:SS), in case the inserted row didn't inherit it.
$ws.Range("D397").NumberFormat = $ws.Range("D398").NumberFormat
